$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '321.51'
    'E2' = '-3.18%'
    'D3' = '42.93'
    'E3' = '-6.44%'
    'E4' = '-8.18%'
    'D5' = '0.08159'
    'E5' = '-2.49%'
    'D6' = '4.327'
    'E6' = '-3.44%'
    'D7' = '1.834'
    'E7' = '-10.96%'
    'D8' = '0.9424'
    'E8' = '-4.31%'
    'D9' = '0.1119'
    'E9' = '-3.49%'
    'D10' = '0.1863'
    'E10' = '-3.20%'
    'D11' = '0.09329'
    'E11' = '-6.40%'
    'D12' = '0.04609'
    'E12' = '-1.26%'
    'D13' = '7.412'
    'E13' = '-28.66%'
    'D14' = '0.1059'
    'E14' = '-0.01%'
    'D15' = '0.001291'
    'E15' = '0.55%'
    'D16' = '0.005696'
    'E16' = '-6.43%'
    'D17' = '3.359'
    'E17' = '-0.46%'
    'D18' = '2.512'
    'E18' = '-2.74%'
    'D19' = '0.3367'
    'E19' = '0.03%'
    'E20' = '-0.89%'
    'E21' = '-1.19%'
    'D22' = '0.04138'
    'E22' = '-1.58%'
    'E23' = '-4.70%'
    'D24' = '0.004299'
    'E24' = '-7.30%'
    'E25' = '-13.49%'
    'D26' = '0.0002981'
    'E26' = '-20.47%'
    'E38' = '-2.61%'
    'D39' = '0.05523'
    'E39' = '-4.89%'
    'D40' = '0.007971'
    'E40' = '3.10%'
    'E41' = '-3.05%'
    'D42' = '0.006552'
    'E42' = '-9.68%'
    'D43' = '0.002120'
    'E43' = '0.02%'
    'D44' = '0.007438'
    'E44' = '-17.77%'
    'E45' = '-6.02%'
    'D46' = '0.00006997'
    'E46' = '-4.07%'
    'D47' = '0.00000000750'
    'E47' = '-0.24%'
    'D48' = '0.003465'
    'E48' = '-0.94%'
    'E49' = '0.76%'
    'D50' = '0.00002100'
    'E50' = '-0.24%'
    'D51' = '0.0002000'
    'E51' = '-0.24%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
